# Add 18 new rows (993-1010) to the "day" sheet with freshly scraped stock data
$wb = $excel.ActiveWorkbook
$wsDay = $wb.Worksheets.Item("day")

$wsDay.Cells.Item(993,1).Value = 1
$wsDay.Cells.Item(993,2).Value = "ULTRACEMCO"
$wsDay.Cells.Item(993,3).Value = "Ultratech Cement Limited"
$wsDay.Cells.Item(993,4).NumberFormat = "@"
$wsDay.Cells.Item(993,4).Value = "532538"
$wsDay.Cells.Item(993,5).Value = 3.98
$wsDay.Cells.Item(993,6).Value = 11648.55
$wsDay.Cells.Item(993,7).Value = 633694
$wsDay.Cells.Item(993,8).Value = "day"
$wsDay.Cells.Item(993,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(994,1).Value = 2
$wsDay.Cells.Item(994,2).Value = "HEROMOTOCO"
$wsDay.Cells.Item(994,3).Value = "Hero Motocorp Limited"
$wsDay.Cells.Item(994,4).NumberFormat = "@"
$wsDay.Cells.Item(994,4).Value = "500182"
$wsDay.Cells.Item(994,5).Value = -0.28
$wsDay.Cells.Item(994,6).Value = 4748.45
$wsDay.Cells.Item(994,7).Value = 510429
$wsDay.Cells.Item(994,8).Value = "day"
$wsDay.Cells.Item(994,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(995,1).Value = 3
$wsDay.Cells.Item(995,2).Value = "HAL"
$wsDay.Cells.Item(995,3).Value = "Hindustan Aeronautics Ltd"
$wsDay.Cells.Item(995,4).NumberFormat = "@"
$wsDay.Cells.Item(995,4).Value = "541154"
$wsDay.Cells.Item(995,5).Value = 0.62
$wsDay.Cells.Item(995,6).Value = 4504.75
$wsDay.Cells.Item(995,7).Value = 1068673
$wsDay.Cells.Item(995,8).Value = "day"
$wsDay.Cells.Item(995,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(996,1).Value = 4
$wsDay.Cells.Item(996,2).Value = "NAVINFLUOR"
$wsDay.Cells.Item(996,3).Value = "Navin Fluorine International Limited"
$wsDay.Cells.Item(996,4).NumberFormat = "@"
$wsDay.Cells.Item(996,4).Value = "532504"
$wsDay.Cells.Item(996,5).Value = 2.14
$wsDay.Cells.Item(996,6).Value = 3583
$wsDay.Cells.Item(996,7).Value = 146051
$wsDay.Cells.Item(996,8).Value = "day"
$wsDay.Cells.Item(996,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(997,1).Value = 5
$wsDay.Cells.Item(997,2).Value = "GODREJPROP"
$wsDay.Cells.Item(997,3).Value = "Godrej Properties Limited"
$wsDay.Cells.Item(997,4).NumberFormat = "@"
$wsDay.Cells.Item(997,4).Value = "533150"
$wsDay.Cells.Item(997,5).Value = 4.51
$wsDay.Cells.Item(997,6).Value = 2901.4
$wsDay.Cells.Item(997,7).Value = 970580
$wsDay.Cells.Item(997,8).Value = "day"
$wsDay.Cells.Item(997,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(998,1).Value = 6
$wsDay.Cells.Item(998,2).Value = "DEEPAKNTR"
$wsDay.Cells.Item(998,3).Value = "Deepak Nitrite Limited"
$wsDay.Cells.Item(998,4).NumberFormat = "@"
$wsDay.Cells.Item(998,4).Value = "506401"
$wsDay.Cells.Item(998,5).Value = 1.17
$wsDay.Cells.Item(998,6).Value = 2759.05
$wsDay.Cells.Item(998,7).Value = 158213
$wsDay.Cells.Item(998,8).Value = "day"
$wsDay.Cells.Item(998,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(999,1).Value = 7
$wsDay.Cells.Item(999,2).Value = "MUTHOOTFIN"
$wsDay.Cells.Item(999,3).Value = "Muthoot Finance Limited"
$wsDay.Cells.Item(999,4).NumberFormat = "@"
$wsDay.Cells.Item(999,4).Value = "533398"
$wsDay.Cells.Item(999,5).Value = 0.97
$wsDay.Cells.Item(999,6).Value = 1935.55
$wsDay.Cells.Item(999,7).Value = 416898
$wsDay.Cells.Item(999,8).Value = "day"
$wsDay.Cells.Item(999,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(1000,1).Value = 8
$wsDay.Cells.Item(1000,2).Value = "DALBHARAT"
$wsDay.Cells.Item(1000,3).Value = "Dalmia Bharat Limited"
$wsDay.Cells.Item(1000,4).NumberFormat = "@"
$wsDay.Cells.Item(1000,4).Value = "533309"
$wsDay.Cells.Item(1000,5).Value = 3.26
$wsDay.Cells.Item(1000,6).Value = 1879.75
$wsDay.Cells.Item(1000,7).Value = 1149388
$wsDay.Cells.Item(1000,8).Value = "day"
$wsDay.Cells.Item(1000,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(1001,1).Value = 9
$wsDay.Cells.Item(1001,2).Value = "HDFCBANK"
$wsDay.Cells.Item(1001,3).Value = "Hdfc Bank Limited"
$wsDay.Cells.Item(1001,4).NumberFormat = "@"
$wsDay.Cells.Item(1001,4).Value = "500180"
$wsDay.Cells.Item(1001,5).Value = 0.48
$wsDay.Cells.Item(1001,6).Value = 1804.7
$wsDay.Cells.Item(1001,7).Value = 7555438
$wsDay.Cells.Item(1001,8).Value = "day"
$wsDay.Cells.Item(1001,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(1002,1).Value = 10
$wsDay.Cells.Item(1002,2).Value = "IPCALAB"
$wsDay.Cells.Item(1002,3).Value = "Ipca Laboratories Limited"
$wsDay.Cells.Item(1002,4).NumberFormat = "@"
$wsDay.Cells.Item(1002,4).Value = "524494"
$wsDay.Cells.Item(1002,5).Value = -0.62
$wsDay.Cells.Item(1002,6).Value = 1533
$wsDay.Cells.Item(1002,7).Value = 328846
$wsDay.Cells.Item(1002,8).Value = "day"
$wsDay.Cells.Item(1002,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(1003,1).Value = 11
$wsDay.Cells.Item(1003,2).Value = "CHOLAFIN"
$wsDay.Cells.Item(1003,3).Value = "Cholamandalam Investment And Finance Company Limited"
$wsDay.Cells.Item(1003,4).NumberFormat = "@"
$wsDay.Cells.Item(1003,4).Value = "511243"
$wsDay.Cells.Item(1003,5).Value = 1.24
$wsDay.Cells.Item(1003,6).Value = 1249.2
$wsDay.Cells.Item(1003,7).Value = 4211650
$wsDay.Cells.Item(1003,8).Value = "day"
$wsDay.Cells.Item(1003,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(1004,1).Value = 12
$wsDay.Cells.Item(1004,2).Value = "GODREJCP"
$wsDay.Cells.Item(1004,3).Value = "Godrej Consumer Products Limited"
$wsDay.Cells.Item(1004,4).NumberFormat = "@"
$wsDay.Cells.Item(1004,4).Value = "532424"
$wsDay.Cells.Item(1004,5).Value = -1.41
$wsDay.Cells.Item(1004,6).Value = 1227.15
$wsDay.Cells.Item(1004,7).Value = 565888
$wsDay.Cells.Item(1004,8).Value = "day"
$wsDay.Cells.Item(1004,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(1005,1).Value = 13
$wsDay.Cells.Item(1005,2).Value = "SUNTV"
$wsDay.Cells.Item(1005,3).Value = "Sun Tv Network Limited"
$wsDay.Cells.Item(1005,4).NumberFormat = "@"
$wsDay.Cells.Item(1005,4).Value = "532733"
$wsDay.Cells.Item(1005,5).Value = -0.18
$wsDay.Cells.Item(1005,6).Value = 757
$wsDay.Cells.Item(1005,7).Value = 262368
$wsDay.Cells.Item(1005,8).Value = "day"
$wsDay.Cells.Item(1005,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(1006,1).Value = 14
$wsDay.Cells.Item(1006,2).Value = "GNFC"
$wsDay.Cells.Item(1006,3).Value = "Gujarat Narmada Valley Fertilizers And Chemicals Limited"
$wsDay.Cells.Item(1006,4).NumberFormat = "@"
$wsDay.Cells.Item(1006,4).Value = "500670"
$wsDay.Cells.Item(1006,5).Value = 2.11
$wsDay.Cells.Item(1006,6).Value = 652.25
$wsDay.Cells.Item(1006,7).Value = 1768312
$wsDay.Cells.Item(1006,8).Value = "day"
$wsDay.Cells.Item(1006,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(1007,1).Value = 15
$wsDay.Cells.Item(1007,2).Value = "APOLLOTYRE"
$wsDay.Cells.Item(1007,3).Value = "Apollo Tyres Limited"
$wsDay.Cells.Item(1007,4).NumberFormat = "@"
$wsDay.Cells.Item(1007,4).Value = "500877"
$wsDay.Cells.Item(1007,5).Value = 0.8
$wsDay.Cells.Item(1007,6).Value = 513.75
$wsDay.Cells.Item(1007,7).Value = 706858
$wsDay.Cells.Item(1007,8).Value = "day"
$wsDay.Cells.Item(1007,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(1008,1).Value = 16
$wsDay.Cells.Item(1008,2).Value = "VEDL"
$wsDay.Cells.Item(1008,3).Value = "Vedanta Limited"
$wsDay.Cells.Item(1008,4).NumberFormat = "@"
$wsDay.Cells.Item(1008,4).Value = "500295"
$wsDay.Cells.Item(1008,5).Value = 1.55
$wsDay.Cells.Item(1008,6).Value = 460.55
$wsDay.Cells.Item(1008,7).Value = 5709578
$wsDay.Cells.Item(1008,8).Value = "day"
$wsDay.Cells.Item(1008,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(1009,1).Value = 17
$wsDay.Cells.Item(1009,2).Value = "COALINDIA"
$wsDay.Cells.Item(1009,3).Value = "Coal India Limited"
$wsDay.Cells.Item(1009,4).NumberFormat = "@"
$wsDay.Cells.Item(1009,4).Value = "533278"
$wsDay.Cells.Item(1009,5).Value = 1.27
$wsDay.Cells.Item(1009,6).Value = 421.7
$wsDay.Cells.Item(1009,7).Value = 6405040
$wsDay.Cells.Item(1009,8).Value = "day"
$wsDay.Cells.Item(1009,9).Value = "02/12/2024 11:35:24"

$wsDay.Cells.Item(1010,1).Value = 18
$wsDay.Cells.Item(1010,2).Value = "CROMPTON"
$wsDay.Cells.Item(1010,3).Value = "Crompton Greaves Consumer Electricals Limited"
$wsDay.Cells.Item(1010,4).NumberFormat = "@"
$wsDay.Cells.Item(1010,4).Value = "539876"
$wsDay.Cells.Item(1010,5).Value = 1.75
$wsDay.Cells.Item(1010,6).Value = 416.85
$wsDay.Cells.Item(1010,7).Value = 1336478
$wsDay.Cells.Item(1010,8).Value = "day"
$wsDay.Cells.Item(1010,9).Value = "02/12/2024 11:35:24"

# Fix bsecode column (D) on the "week" sheet: these 26 rows were stored as text
# instead of numbers; convert them to proper numeric values.
$wsWeek = $wb.Worksheets.Item("week")
$wsWeek.Cells.Item(613,4).Value = 532466
$wsWeek.Cells.Item(614,4).Value = 532541
$wsWeek.Cells.Item(615,4).Value = 505200
$wsWeek.Cells.Item(616,4).Value = 540762
$wsWeek.Cells.Item(617,4).Value = 500495
$wsWeek.Cells.Item(618,4).Value = 532175
$wsWeek.Cells.Item(619,4).Value = 500271
$wsWeek.Cells.Item(620,4).Value = 500770
$wsWeek.Cells.Item(621,4).Value = 543220
$wsWeek.Cells.Item(622,4).Value = 539268
$wsWeek.Cells.Item(623,4).Value = 532508
$wsWeek.Cells.Item(624,4).Value = 543300
$wsWeek.Cells.Item(625,4).Value = 540777
$wsWeek.Cells.Item(626,4).Value = 500253
$wsWeek.Cells.Item(627,4).Value = 532814
$wsWeek.Cells.Item(628,4).Value = 532810
$wsWeek.Cells.Item(629,4).Value = 500049
$wsWeek.Cells.Item(630,4).Value = 500103
$wsWeek.Cells.Item(631,4).Value = 532210
$wsWeek.Cells.Item(632,4).Value = 531213
$wsWeek.Cells.Item(633,4).Value = 543257
$wsWeek.Cells.Item(634,4).Value = 500183
$wsWeek.Cells.Item(635,4).Value = 532477
$wsWeek.Cells.Item(636,4).Value = 532149
$wsWeek.Cells.Item(637,4).Value = 532461
$wsWeek.Cells.Item(638,4).Value = 533098
